$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(19, 44738, "AM", "Maggie", "Positive"),
    @(20, 44740, "AM", "Maggie", "Negative"),
    @(21, 44740, "AM", "Michael", "Positive"),
    @(22, 44740, "AM", "Jessica", "Positive"),
    @(23, 44736, "PM", "Michael", "Negative"),
    @(24, 44737, "AM", "Michael", "Positive"),
    @(25, 44742, "AM", "Michael", "Positive"),
    @(26, 44744, "AM", "Michael", "Positive"),
    @(27, 44742, "AM", "Jessica", "Positive"),
    @(28, 44737, "AM", "Ben", "Positive")
)

foreach ($row in $data) {
    $r = $row[0]
    $dateSerial = $row[1]
    $time = $row[2]
    $name = $row[3]
    $test = $row[4]

    $ws.Cells.Item($r, 1).Value = [DateTime]::FromOADate($dateSerial)
    $ws.Cells.Item($r, 2).NumberFormat = "d-mmm"
    $ws.Cells.Item($r, 2).Value = $time
    $ws.Cells.Item($r, 3).Value = $name
    $ws.Cells.Item($r, 4).Value = $test
}

$ws.Range("D28").Select()
